# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E) for worker LUIS GUSTAVO MORENO OLIVEROS
# (rows 19-23) was entered in descending order (1912, 1911, 1910, 1909,
# 1908). Correct it to ascending chronological order (1908, 1909, 1910,
# 1911, 1912). Row 21 (1910) already sits in the middle and keeps its
# value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E19").Value = "1908"
$ws.Range("E20").Value = "1909"
$ws.Range("E21").Value = "1910"
$ws.Range("E22").Value = "1911"
$ws.Range("E23").Value = "1912"
